# DP: create_forecast_basic current - palestinian_from_demo_230622.ipynb - change TAZ_V4 input date
# Adds a new "TAZ_V4_date" input row (row 7) below the existing inputs/outputs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label/value pair for the TAZ_V4 input date.
$ws.Range("A7").Value = "TAZ_V4_date"
$ws.Range("B7").Value = 240404

# Move the active selection to the newly added row, mirroring what a user
# would see right after typing the new entry in (A7:B7, anchored at A7).
$excel.Goto($ws.Range("A7:B7"))
